$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update rows 3 through 91 with the new/reshuffled localization content.
# (Rows 1-2 are unchanged headers/welcome row.)
$ws.Cells.Item(3, 1).Value = "title"
$ws.Cells.Item(3, 2).Value = "<size=50>Attack on Blob</size>\nMultiply and Divide"
$ws.Cells.Item(3, 2).WrapText = $true

$ws.Cells.Item(4, 1).Value = "credits"
$ws.Cells.Item(4, 2).Value = "CREDITS"
$ws.Cells.Item(4, 2).WrapText = $true

$ws.Cells.Item(5, 1).Value = "credits_desc"
$ws.Cells.Item(5, 2).Value = "Written by: David Dionisio\nMusic from: Kevin Macleod"
$ws.Cells.Item(5, 2).WrapText = $true

$ws.Cells.Item(6, 1).Value = "play"
$ws.Cells.Item(6, 2).Value = "PLAY"

$ws.Cells.Item(7, 1).Value = "multiply"
$ws.Cells.Item(7, 2).Value = "Multiply"

$ws.Cells.Item(8, 1).Value = "divide"
$ws.Cells.Item(8, 2).Value = "Divide"

$ws.Cells.Item(9, 1).Value = "ready"
$ws.Cells.Item(9, 2).Value = "READY"

$ws.Cells.Item(10, 1).Value = "set"
$ws.Cells.Item(10, 2).Value = "SET"

$ws.Cells.Item(11, 1).Value = "go"
$ws.Cells.Item(11, 2).Value = "GO"

$ws.Cells.Item(12, 1).Value = "victory"
$ws.Cells.Item(12, 2).Value = "VICTORY"

$ws.Cells.Item(13, 1).Value = "score"
$ws.Cells.Item(13, 2).Value = "SCORE"

$ws.Cells.Item(14, 1).Value = "time"
$ws.Cells.Item(14, 2).Value = "TIME"

$ws.Cells.Item(15, 1).Value = "time_bonus"
$ws.Cells.Item(15, 2).Value = "TIME BONUS"

$ws.Cells.Item(16, 1).Value = "perfect"
$ws.Cells.Item(16, 2).Value = "PERFECT"

$ws.Cells.Item(17, 1).Value = "total"
$ws.Cells.Item(17, 2).Value = "TOTAL"

$ws.Cells.Item(18, 1).Value = "options"
$ws.Cells.Item(18, 2).Value = "OPTIONS"

$ws.Cells.Item(19, 1).Value = "music"
$ws.Cells.Item(19, 2).Value = "MUSIC"

$ws.Cells.Item(20, 1).Value = "sound"
$ws.Cells.Item(20, 2).Value = "SOUND"

$ws.Cells.Item(21, 1).Value = "speech"
$ws.Cells.Item(21, 2).Value = "SPEECH"

$ws.Cells.Item(22, 1).Value = "on"
$ws.Cells.Item(22, 2).Value = "ON"

$ws.Cells.Item(23, 1).Value = "off"
$ws.Cells.Item(23, 2).Value = "OFF"

$ws.Cells.Item(24, 1).Value = "close"
$ws.Cells.Item(24, 2).Value = "CLOSE"

$ws.Cells.Item(25, 1).Value = "complete"
$ws.Cells.Item(25, 2).Value = "COMPLETE"

$ws.Cells.Item(26, 1).Value = "mult2_title"
$ws.Cells.Item(26, 2).Value = "Multiples of 2"

$ws.Cells.Item(27, 1).Value = "mult3_title"
$ws.Cells.Item(27, 2).Value = "Multiples of 3"

$ws.Cells.Item(28, 1).Value = "mult4_title"
$ws.Cells.Item(28, 2).Value = "Multiples of 4"

$ws.Cells.Item(29, 1).Value = "mult5_title"
$ws.Cells.Item(29, 2).Value = "Multiples of 5"

$ws.Cells.Item(30, 1).Value = "mult6_title"
$ws.Cells.Item(30, 2).Value = "Multiples of 6"

$ws.Cells.Item(31, 1).Value = "mult7_title"
$ws.Cells.Item(31, 2).Value = "Multiples of 7"

$ws.Cells.Item(32, 1).Value = "mult8_title"
$ws.Cells.Item(32, 2).Value = "Multiples of 8"

$ws.Cells.Item(33, 1).Value = "mult9_title"
$ws.Cells.Item(33, 2).Value = "Multiples of 9"

$ws.Cells.Item(34, 1).Value = "mult10_title"
$ws.Cells.Item(34, 2).Value = "Multiples of 10"

$ws.Cells.Item(35, 1).Value = "multiplicationTable"
$ws.Cells.Item(35, 2).Value = "Multiplication Table"

$ws.Cells.Item(36, 1).Value = "mult10_other_title"
$ws.Cells.Item(36, 2).Value = "...also 100, 1000, and so forth."

$ws.Cells.Item(37, 1).Value = "multTable_instruct"
$ws.Cells.Item(37, 2).Value = "Press this button to review the multiplication table."
$ws.Cells.Item(37, 3).Value = 5

$ws.Cells.Item(38, 1).Value = "proceed_instruct"
$ws.Cells.Item(38, 2).Value = "Press this button to proceed."
$ws.Cells.Item(38, 3).Value = 5

$ws.Cells.Item(39, 1).Value = "multiplication"
$ws.Cells.Item(39, 2).Value = "Multiplication"

$ws.Cells.Item(40, 1).Value = "division"
$ws.Cells.Item(40, 2).Value = "Division"

$ws.Cells.Item(41, 1).Value = "example"
$ws.Cells.Item(41, 2).Value = "Example"

$ws.Cells.Item(42, 1).Value = "commutative_title"
$ws.Cells.Item(42, 2).Value = "Commutative Property"

$ws.Cells.Item(43, 1).Value = "associative_title"
$ws.Cells.Item(43, 2).Value = "Associative Property"

$ws.Cells.Item(44, 1).Value = "distributive_title"
$ws.Cells.Item(44, 2).Value = "Distributive Property"

$ws.Cells.Item(45, 1).Value = "not_commutative"
$ws.Cells.Item(45, 2).Value = "Not Commutative!"

$ws.Cells.Item(46, 1).Value = "not_associative"
$ws.Cells.Item(46, 2).Value = "Not Associative!"

$ws.Cells.Item(47, 1).Value = "lesson_1_intro_1"
$ws.Cells.Item(47, 2).Value = "Before we proceed, let's first learn some tricks with multiplication!"

$ws.Cells.Item(48, 1).Value = "lesson_1_mult2_1"
$ws.Cells.Item(48, 2).Value = "In multiples of two, the trick is to simply double the number."

$ws.Cells.Item(49, 1).Value = "lesson_1_mult2_2"
$ws.Cells.Item(49, 2).Value = "For example: 2 x 6 can be 6 + 6, which equals to 12."

$ws.Cells.Item(50, 1).Value = "lesson_1_commutative_1"
$ws.Cells.Item(50, 2).Value = "The commutative property means that multiplying numbers in any order gives the same answer."

$ws.Cells.Item(51, 1).Value = "lesson_1_commutative_2"
$ws.Cells.Item(51, 2).Value = "For example: 2 x 3 and 3 x 2 equal 6."

$ws.Cells.Item(52, 1).Value = "lesson_1_commutative_3"
$ws.Cells.Item(52, 2).Value = "With this trick, you only have to remember half the multiplication table!"

$ws.Cells.Item(53, 1).Value = "lesson_1_tutorial_1"
$ws.Cells.Item(53, 2).Value = "Now banish these blobs by connecting them in the correct order using multiplication."

$ws.Cells.Item(54, 1).Value = "lesson_1_tutorial_end_1"
$ws.Cells.Item(54, 2).Value = "Excellent! You are now ready for the mission!"

$ws.Cells.Item(55, 1).Value = "lesson_2_intro_1"
$ws.Cells.Item(55, 2).Value = "Good work! Now it's time to step up the game with multiples of 3 and 4, along with division."

$ws.Cells.Item(56, 1).Value = "lesson_2_intro_2"
$ws.Cells.Item(56, 2).Value = "Let me show you some neat tricks."

$ws.Cells.Item(57, 1).Value = "lesson_2_mult3_1"
$ws.Cells.Item(57, 2).Value = "In multiples of three: double the number, and then add the original number."

$ws.Cells.Item(58, 1).Value = "lesson_2_mult3_2"
$ws.Cells.Item(58, 2).Value = "For example, 3 x 6: double 6 to get 12, and then add 6 to get 18."

$ws.Cells.Item(59, 1).Value = "lesson_2_mult4_1"
$ws.Cells.Item(59, 2).Value = "In multiples of four: double the number, and then double it again."

$ws.Cells.Item(60, 1).Value = "lesson_2_mult4_2"
$ws.Cells.Item(60, 2).Value = "For example, 4 x 6: double 6 to get 12, and then double 12 to get 24."

$ws.Cells.Item(61, 1).Value = "lesson_2_div_1"
$ws.Cells.Item(61, 2).Value = "When it comes to division, think of it as the opposite of multiplication."

$ws.Cells.Item(62, 1).Value = "lesson_2_div_2"
$ws.Cells.Item(62, 2).Value = "Rearranging the equation, and replacing division with multiplication can help."

$ws.Cells.Item(63, 1).Value = "lesson_2_div_3"
$ws.Cells.Item(63, 2).Value = "Unlike multiplication, division is not commutative. So the order of the numbers cannot be changed."

$ws.Cells.Item(64, 1).Value = "lesson_2_tutorial_1"
$ws.Cells.Item(64, 2).Value = "For the next mission, some blobs must be matched with division. Go ahead and try it out."

$ws.Cells.Item(65, 1).Value = "lesson_2_tutorial_end_1"
$ws.Cells.Item(65, 2).Value = "Excellent! You are now ready for the mission!"

$ws.Cells.Item(66, 1).Value = "lesson_3_intro_1"
$ws.Cells.Item(66, 2).Value = "So far so good! Now it's time to go over the multiples of 5 and 6."

$ws.Cells.Item(67, 1).Value = "lesson_3_mult5_1"
$ws.Cells.Item(67, 2).Value = "There are several ways to go about with the multiples of 5."

$ws.Cells.Item(68, 1).Value = "lesson_3_mult5_2"
$ws.Cells.Item(68, 2).Value = "You can count in fives, alternating between 0's and 5's."

$ws.Cells.Item(69, 1).Value = "lesson_3_mult5_3"
$ws.Cells.Item(69, 2).Value = "Or, you can multiply the number by 10, and then half it."

$ws.Cells.Item(70, 1).Value = "lesson_3_mult6_1"
$ws.Cells.Item(70, 2).Value = "In multiples of 6, a good approach is to first multiply the number by 5, and then add the original number."

$ws.Cells.Item(71, 1).Value = "lesson_3_mult6_2"
$ws.Cells.Item(71, 2).Value = "For example, 6 x 7: multiply 7 by 5 to get 35, and then add 7 to get 42."

$ws.Cells.Item(72, 1).Value = "lesson_3_associative_1"
$ws.Cells.Item(72, 2).Value = "The associative property means that we can group any of the numbers in the equation."

$ws.Cells.Item(73, 1).Value = "lesson_3_associative_2"
$ws.Cells.Item(73, 2).Value = "Numbers that are grouped together will be calculated first."

$ws.Cells.Item(74, 1).Value = "lesson_3_associative_3"
$ws.Cells.Item(74, 2).Value = "This trick is useful in multiplication if you want to split up a large number into several multiples."

$ws.Cells.Item(75, 1).Value = "lesson_3_associative_4"
$ws.Cells.Item(75, 2).Value = "Remember that unlike multiplication, division is not associative."

$ws.Cells.Item(76, 1).Value = "lesson_3_end_1"
$ws.Cells.Item(76, 2).Value = "Now let us proceed to the next mission!"

$ws.Cells.Item(77, 1).Value = "lesson_4_intro_1"
$ws.Cells.Item(77, 2).Value = "We are almost there! Let's take a look at some tricks for the multiples of 7 and 8."

$ws.Cells.Item(78, 1).Value = "lesson_4_mult7_1"
$ws.Cells.Item(78, 2).Value = "In multiples of 7, a good trick is to multiply the number by 5, and then add the original number twice."

$ws.Cells.Item(79, 1).Value = "lesson_4_mult7_2"
$ws.Cells.Item(79, 2).Value = "You can also use the multiples of 6: multiply the number by 6, and then add the original number."

$ws.Cells.Item(80, 1).Value = "lesson_4_mult8_1"
$ws.Cells.Item(80, 2).Value = "In multiples of 8: double the number three times."

$ws.Cells.Item(81, 1).Value = "lesson_4_mult8_2"
$ws.Cells.Item(81, 2).Value = "For example, 8 x 4: double 4 to get 8, double 8 to get 16, and finally double 16 to get 32."

$ws.Cells.Item(82, 1).Value = "lesson_4_distributive_1"
$ws.Cells.Item(82, 2).Value = "The distributive property allows you to distribute a number to a group."

$ws.Cells.Item(83, 1).Value = "lesson_4_distributive_2"
$ws.Cells.Item(83, 2).Value = "For multiplication, this is a good trick for splitting up a number, and computing each one separately."

$ws.Cells.Item(84, 1).Value = "lesson_4_distributive_3"
$ws.Cells.Item(84, 2).Value = "As they say: when an obstacle is too large to handle, divide and conquer."

$ws.Cells.Item(85, 1).Value = "lesson_4_end_1"
$ws.Cells.Item(85, 2).Value = "Now onwards to the next mission!"

$ws.Cells.Item(86, 1).Value = "lesson_5_intro_1"
$ws.Cells.Item(86, 2).Value = "We are just one mission away from total victory! Let's look at the last two multiples: 9 and 10."

$ws.Cells.Item(87, 1).Value = "lesson_5_mult9_1"
$ws.Cells.Item(87, 2).Value = "In multiples of 9: simply multiply the number by 10, and then subtract the original number."

$ws.Cells.Item(88, 1).Value = "lesson_5_mult9_2"
$ws.Cells.Item(88, 2).Value = "For example, 9 x 6: multiply 6 by 10 to get 60, and then subtract 6 to get 54."

$ws.Cells.Item(89, 1).Value = "lesson_5_mult10_1"
$ws.Cells.Item(89, 2).Value = "In multiples of 10: just add a 0 at the end of the number."

$ws.Cells.Item(90, 1).Value = "lesson_5_mult10_2"
$ws.Cells.Item(90, 2).Value = "In fact, you can do this for any amount of 0's such as: 100, 1000, etc."

$ws.Cells.Item(91, 1).Value = "lesson_5_end_1"
$ws.Cells.Item(91, 2).Value = "Now it's time to clean up the last remaining blobs, good luck!"

# Clear stale numeric values left over in column C from rows that moved.
$ws.Cells.Item(33, 3).ClearContents()
$ws.Cells.Item(34, 3).ClearContents()

# Update the selected/active cell shown when the workbook is opened.
$ws.Range("B4").Select()

